$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 301
$ws.Range("F5").Value = 5338
$ws.Range("F6").Value = 617
$ws.Range("F7").Value = 454
$ws.Range("F8").Value = 243
$ws.Range("F9").Value = 1106
$ws.Range("F11").Value = 158
$ws.Range("F12").Value = 56
$ws.Range("F13").Value = 748
$ws.Range("F14").Value = 382
$ws.Range("F17").Value = 191
$ws.Range("F19").Value = 375
$ws.Range("F20").Value = 6158
$ws.Range("F22").Value = 49
$ws.Range("F24").Value = 7123
$ws.Range("F27").Value = 3272
$ws.Range("F28").Value = 390
$ws.Range("F29").Value = 779
$ws.Range("F31").Value = 327
$ws.Range("F32").Value = 149
$ws.Range("F33").Value = 154
$ws.Range("F34").Value = 1208
$ws.Range("F35").Value = 113
$ws.Range("F36").Value = 33
$ws.Range("F38").Value = 960
$ws.Range("F39").Value = 1214

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 54
$ws.Range("F5").Value = 68

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1160

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1160
$ws.Range("F7").Value = 301
$ws.Range("F8").Value = 5338
$ws.Range("F9").Value = 617
$ws.Range("F10").Value = 454
$ws.Range("F11").Value = 243
$ws.Range("F12").Value = 1106
$ws.Range("F14").Value = 158
$ws.Range("F15").Value = 56
$ws.Range("F16").Value = 748
$ws.Range("F17").Value = 382
$ws.Range("F19").Value = 54
$ws.Range("F21").Value = 191
$ws.Range("F23").Value = 375
$ws.Range("F24").Value = 6158
$ws.Range("F26").Value = 49
$ws.Range("F28").Value = 7123
$ws.Range("F31").Value = 3272
$ws.Range("F32").Value = 390
$ws.Range("F33").Value = 779
$ws.Range("F35").Value = 327
$ws.Range("F37").Value = 149
$ws.Range("F38").Value = 154
$ws.Range("F39").Value = 1208
$ws.Range("F40").Value = 113
$ws.Range("F41").Value = 33
$ws.Range("F43").Value = 960
$ws.Range("F44").Value = 1214
$ws.Range("F48").Value = 68
